## !HACKERRANK - PROGRESS.xlsx -- "Add files via upload"
## Renames Sheet1 -> Progress, adds a new "Plan" sheet with a short list of
## next HackerRank problems to tackle, and records the latest progress-table
## entry (row 15) on the Progress sheet. Also clears the leftover fill on
## the most recent earlier row (16) now that it is no longer the "latest"
## highlighted row.

$wb = $excel.ActiveWorkbook
$progress = $wb.Worksheets.Item(1)
$progress.Name = "Progress"

# ---- Progress sheet: fill in the newest tracked entry (row 15) ----------
$progress.Range("N15").Value = 45126
$progress.Range("O15").Value = "1367.97/2200"
$progress.Range("P15").Value = 96935
$progress.Range("Q15").Formula = "=IF(ROW()>2,(`$P`$2-P15)/`$P`$2,""NA"")"

# Row 16 no longer needs the (stray) fill that had been applied earlier -
# clear it back to "No Fill" so it matches the rest of the table again.
$progress.Range("B16:E16").Interior.Pattern = -4142

# ---- New "Plan" sheet: next problems to solve ----------------------------
$plan = $wb.Worksheets.Add($null, $progress)
$plan.Name = "Plan"

$plan.Range("A2").Value = "Organizing Containers of Balls"
$plan.Range("A3").Value = "Encryption"
$plan.Range("A4").Value = "Bigger is Greater"
$plan.Range("A5").Value = "The Time in Words"
$plan.Range("A1").Value = "Plan - next problems to solve"
$plan.Range("A2").Select() | Out-Null

# Leave the selection where the author last clicked before saving.
$progress.Activate()
$progress.Range("M16").Select() | Out-Null
